# Regenerate the handoff report: a new handoff run produced a new request
# id (b9b62266-5f5a-4ea8-a245-52d179bcbf78) and new target-file hashes,
# replacing the previous id (0e175488-7063-4752-a540-204b40e961cb) and
# bumping the "latest handoff" timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$oldId = "0e175488-7063-4752-a540-204b40e961cb"
$newId = "b9b62266-5f5a-4ea8-a245-52d179bcbf78"

$oldZhHash = "6834214d81b3a79956ad596d9b14c7efd7529b8f"
$newZhHash = "acf5a99546d2ece5fa2a38219468f29fdbf9c500"
$oldDeHash = "6834214d81b3a79956ad596d9b14c7efd7529b8f"
$newDeHash = "acf5a99546d2ece5fa2a38219468f29fdbf9c500"

$mdName = $newId + ".md"
$zhXlfName = $newId + "." + $newZhHash + ".zh-cn.xlf"
$deXlfName = $newId + "." + $newDeHash + ".de-de.xlf"

# NOTE: only the displayed link text / cell text changed upstream - the
# underlying hyperlink targets (relationship URLs) stay exactly as they
# were, so we re-add each hyperlink against its ORIGINAL (old-id) address.
$mdAddr = "https://github.com/OpenLocalizationTest/oltest/blob/87660e3233d7841ac49b57b98ad76d70bfd99ddc/e2e/" + $oldId + ".md"
$zhXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74424c9fa392ac0f265ce1ecdd1de938c167858a/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/" + $oldId + "." + $oldZhHash + ".zh-cn.xlf"
$deXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e030dbd89ac580c9f5c2cf7d95358ba1413eb65/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/" + $oldId + "." + $oldDeHash + ".de-de.xlf"

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = $mdName
$ws.Range("D2").Value = "2016-03-20 17:23:33"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    $mdAddr,
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$ws.Range("A2").Font.Underline = 2
$ws.Range("A2").Font.Color = 15570276

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = $mdName
$ws.Range("D2").Value = $zhXlfName
$ws.Range("E2").Value = "2016-03-20 17:23:25"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    $mdAddr,
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$ws.Range("A2").Font.Underline = 2
$ws.Range("A2").Font.Color = 15570276

$ws.Hyperlinks.Add(
    $ws.Range("D2"),
    $zhXlfAddr,
    [Type]::Missing,
    [Type]::Missing,
    $zhXlfName) | Out-Null
$ws.Range("D2").Font.Underline = 2
$ws.Range("D2").Font.Color = 15570276

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = $mdName
$ws.Range("D2").Value = $deXlfName
$ws.Range("E2").Value = "2016-03-20 17:23:33"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    $mdAddr,
    [Type]::Missing,
    [Type]::Missing,
    $mdName) | Out-Null
$ws.Range("A2").Font.Underline = 2
$ws.Range("A2").Font.Color = 15570276

$ws.Hyperlinks.Add(
    $ws.Range("D2"),
    $deXlfAddr,
    [Type]::Missing,
    [Type]::Missing,
    $deXlfName) | Out-Null
$ws.Range("D2").Font.Underline = 2
$ws.Range("D2").Font.Color = 15570276

Write-Host "Handoff report regenerated"
